$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: DATE_TYPE_CODE changes from "001" to "002" (keep as text, avoid numeric coercion)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"

# N2: REPORT_DATE changes to the new period end date (text, same format as existing cells)
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Updated numeric figures for the new reporting period
$ws.Range("O2").Value = 550799516.5599999
$ws.Range("P2").Value = 7562468.16
$ws.Range("Q2").Value = 202844124.41
$ws.Range("S2").Value = 247143610.93
$ws.Range("U2").Value = 2713684.86
$ws.Range("W2").Value = 235709914.77
$ws.Range("X2").Value = 60553904.24
$ws.Range("AB2").Value = 315089601.79
$ws.Range("AF2").Value = 206.3139881775
$ws.Range("AG2").Value = 42.7941397338

# These ratio columns no longer carry a value in the updated data (now blank)
$ws.Range("R2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
